# Rename the paired "old/new" comparison headers to the format-version
# specific names (e.g. "Segmentname_old" -> "Segmentname_FV2210",
# "Segmentname_new" -> "Segmentname_FV2304"), freeze the header row, and
# wrap the data range in a native Excel Table so filtering/sorting keeps
# the new header names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseCols = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")
$leftLetters = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$rightLetters = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

for ($i = 0; $i -lt $baseCols.Length; $i++) {
    $ws.Range($leftLetters[$i] + "1").Value = $baseCols[$i] + "_FV2210"
    $ws.Range($rightLetters[$i] + "1").Value = $baseCols[$i] + "_FV2304"
}

# Freeze the header row (row 1) so it stays visible while scrolling.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the data range into an Excel Table (ListObject) using the renamed headers.
$rng = $ws.Range("A1:U64")
$tbl = $ws.ListObjects.Add(1, $rng, [System.Type]::Missing, 1)
$tbl.Name = "Table1"
